$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Age"
$ws.Range("C1").Value = "Sex"
$ws.Range("D1").Value = "Weight kg"
$ws.Range("E1").Value = "Height cm"
$ws.Range("F1").Value = "Weight lb"
$ws.Range("G1").Value = "Height ft'in`""
$ws.Range("H1").Value = "Activity Level"
$ws.Range("I1").Value = "Goal"
